$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "MEC-3B-Fresagem"
$ws.Range("C3").Value = "MEC-3B-Fresagem"
$ws.Range("D4").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("D7").Value = "MEC-3B-Fresagem"
$ws.Range("D8").Value = "MEC-3B-Fresagem"
$ws.Range("E10").Value = "MEC-3A-M. A. Comp; Cad / CAM"
$ws.Range("E11").Value = "MEC-3A-Fresagem"
$ws.Range("C12").Value = "-"
$ws.Range("E12").Value = "MEC-3A-Fresagem"
$ws.Range("C14").Value = "MEC-3A-Fresagem"
$ws.Range("E14").Value = "MEC-3A-Fresagem"
$ws.Range("C15").Value = "-"
$ws.Range("E15").Value = "MEC-3A-M. A. Comp; Cad / CAM"
$ws.Range("B16").Value = "MEC-3A-M. A. Comp; Cad / CAM"
$ws.Range("E16").Value = "MEC-3A-M. A. Comp; Cad / CAM"
